$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing-cell value fixes ---
$ws.Range("Q57").Value = 0
$ws.Range("R272").Value = 0
$ws.Range("O273").Value = 1
$ws.Range("R273").Value = 0

# --- New weekly rows 274-278 (continuing the OHLCV/date-part table) ---
$newRows = @(
    @(274, 45474, 6750,            6786,          6544.10009765625, 6634.10009765625, 6603.69921875,   2533461, 2024, 7, 1,  0, 0, 0, 27, 0, 0, 0),
    @(275, 45481, 6639,            6675,          6311.10009765625, 6581,             6550.84228515625, 2962781, 2024, 7, 8,  0, 0, 0, 28, 0, 0, 0),
    @(276, 45488, 6591,            6717.85009765625, 6257.5,       6349.85009765625, 6349.85009765625, 3202220, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
    @(277, 45495, 6330,            6564.75,       5930.0498046875,  6548.5,           6548.5,           3003488, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(278, 45502, 6584.9501953125, 6882,          6540,             6858.2001953125,  6858.2001953125,  1724736, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
    # Column R (18) stays blank for these new rows, matching the source data.
}
